$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1346003333333333
$ws.Range("H2").Value = 0.403801
$ws.Range("I2").Value = 0.009651054304565105
$ws.Range("J2").Value = 0.009651054304565105
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 0.3666515323338889
$ws.Range("R2").Value = 3.299863791005
$ws.Range("S2").Value = 0.0004473746680569466
$ws.Range("T2").Value = 0.0004473746680569466
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1346003333333333
$ws.Range("H3").Value = 0.403801
$ws.Range("I3").Value = 0.009651054304565105
$ws.Range("J3").Value = 0.009651054304565105
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 5.467550427943557
$ws.Range("R3").Value = 49.20795385149201
$ws.Range("S3").Value = 0.006671303245934319
$ws.Range("T3").Value = 0.006671303245934318
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1346003333333333
$ws.Range("H4").Value = 0.403801
$ws.Range("I4").Value = 0.009651054304565105
$ws.Range("J4").Value = 0.009651054304565105
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 2.075440900761667
$ws.Range("R4").Value = 18.678968106855
$ws.Range("S4").Value = 0.002532376390573839
$ws.Range("T4").Value = 0.002532376390573839
$ws.Range("I5").Value = 0.8124788779145131
$ws.Range("J5").Value = 0.8124788779145132
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 30.86674431366167
$ws.Range("R5").Value = 277.800698822955
$ws.Range("S5").Value = 0.0376624622387994
$ws.Range("T5").Value = 0.03766246223879941
$ws.Range("I6").Value = 0.8124788779145131
$ws.Range("J6").Value = 0.8124788779145132
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.561627031040565
$ws.Range("T6").Value = 0.561627031040565
$ws.Range("I7").Value = 0.8124788779145131
$ws.Range("J7").Value = 0.8124788779145132
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 174.722039791145
$ws.Range("R7").Value = 1572.498358120305
$ws.Range("S7").Value = 0.2131893846351487
$ws.Range("T7").Value = 0.2131893846351487
$ws.Range("G8").Value = 2.4807
$ws.Range("H8").Value = 7.4421
$ws.Range("I8").Value = 0.1778700677809217
$ws.Range("J8").Value = 0.1778700677809217
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 6.7574309345
$ws.Range("R8").Value = 60.8168784105
$ws.Range("S8").Value = 0.008245167835509576
$ws.Range("T8").Value = 0.008245167835509576
$ws.Range("G9").Value = 2.4807
$ws.Range("H9").Value = 7.4421
$ws.Range("I9").Value = 0.1778700677809217
$ws.Range("J9").Value = 0.1778700677809217
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 100.7675984948
$ws.Range("R9").Value = 906.9083864532
$ws.Range("S9").Value = 0.1229529047391359
$ws.Range("T9").Value = 0.1229529047391358
$ws.Range("G10").Value = 2.4807
$ws.Range("H10").Value = 7.4421
$ws.Range("I10").Value = 0.1778700677809217
$ws.Range("J10").Value = 0.1778700677809217
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 38.2506202995
$ws.Range("R10").Value = 344.2555826955
$ws.Range("S10").Value = 0.04667199520627628
$ws.Range("T10").Value = 0.04667199520627628
